$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H ("DESCRIPCION DOCUMENTO"), shifting old H:O to I:P
$ws.Columns("H").Insert()
$ws.Columns("H").ColumnWidth = 25.5

# Header for the new column
$ws.Range("H1").Value = "DESCRIPCION DOCUMENTO"

# Sample/test description value typed in as text (with leading apostrophe -> quote-prefixed text)
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "'Descripcion de prueba"

# Copy the formatting of H2 down to H6 (next group header row), no value
$ws.Range("H2").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the view: zoom level and active selection
$ws.Range("G12").Select()
$excel.ActiveWindow.Zoom = 112
